$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exposure conditions")

# Header for new column N (copy formatting from M1, the last header cell)
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "PrecisionTox short identifier"

# New identifiers for rows 2-25 (data rows)
$ids = @(
    "XAA---LA1",
    "XAA---LA2",
    "XAA---LA3",
    "XAA---LA4",
    "XAA---LA1",
    "XAA---LA2",
    "XAA---LA3",
    "XAA---LA4",
    "XAA---LB1",
    "XAA---LB2",
    "XAA---LB3",
    "XAA---LB4",
    "XAA---LB1",
    "XAA---LB2",
    "XAA---LB3",
    "XAA---LB4",
    "XAA---LC1",
    "XAA---LC2",
    "XAA---LC3",
    "XAA---LC4",
    "XAA---LC1",
    "XAA---LC2",
    "XAA---LC3",
    "XAA---LC4"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 14).Value = $ids[$i]
}

# Rows 26-27: update replicate (column I) values and add identifiers
$ws.Range("I26").Value = 1
$ws.Range("I27").Value = 2
$ws.Range("N26").Value = "XAA998ZS1"
$ws.Range("N27").Value = "XAA998ZS2"

$wb.Save()
